$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend the bug-report text in C40 ("Ajout d'un retard de 0 minutes" test)
#    with the extra remark about converting the two numbers into real dates.
#    Editing the cell text (rather than the shared string in place) is required
#    because the original text is shared with C35, which must stay unchanged.
$oldText = $ws.Range("C40").Value()
$ws.Range("C40").Value = $oldText + " + au lieu de concaténer les deux nombres, les convertir directement en dates"

# 2) Insert a new blank row above the "Historique absences/retards" section
#    header (old row 47), pushing it and everything below down by one row.
$ws.Rows("47:47").Insert()

# 3) Restore the view state: select C41, matching where the author was
#    working when they saved the file (the headless runtime has no Window
#    object model, so the scrolled topLeftCell position cannot be set here).
$ws.Range("C41").Select()
